$d = $word.ActiveDocument

# --- Part 1: merge the author-line runs into a single run ---
$rng = $d.Content
$null = $rng.Find.Execute("B. Comnes and A. La Rosa", $false, $false, $false, $false, $false, $true, 1, $false, "B. Comnes and A. La Rosa", 2)

# --- Part 2: insert the lab-instructions paragraphs ---
$anchorPara = $d.Paragraphs.Last
$insertRange = $anchorPara.Range
$insertRange.Collapse(1)

$fragment = @'
<w:p><w:r><w:t>Open the ISE project Navigator</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>File - &gt; New Project</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Give your project a name with no spaces or special characters.  DON’T USE hyphens! (-)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Under top-level source choose HDL</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Click next</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Under evaluation development board choose Spartan 6 sp601 evaluation platform</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Product category family speed package will all be chosen for you</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Under Preferred language choose VHDL</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Click next</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>You will get a project summary window (new project 4.png)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">This closes the new project wizard, and dropes you into your ISE project navigator Integrated Development enviroment.  It is a very non-intuitive interface, that you will only become accustom to with experience.  </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>We will be building a simple And circuit using VHDL and manual schematic layouts.  Both are achieve the same thing, although VHDL will eventually get you more mileage.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Add a new source schematic to your project by selecting Project -&gt; New Source</w:t></w:r><w:r><w:t xml:space="preserve"> which oopens the new source wizard (new source wizard 1)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Give it a name, avoiding special characters and click next to get to the summary of the wizzarrd ( new source wizard 2)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Click finished</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>It will create the schematic and open it. (new schematic 1)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Click the add symbol toolbar button in the vertical toolbar strip</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>This opens a menu to the left of the toolbar with digital electronic components.  They are sorted into categories for somewhat quick access, but the name filter is the quickest way to find what you are looking for if you know the name</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr><w:r><w:t>Click the logic category, and notice the symbols list shortened to only show the logic componants in alphabetical oder</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr><w:r><w:t>Choose the and2 symbole from the symbol menu.</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr><w:r><w:t>Now when you roll over the schematic view, you get something that looks like that compotant following around a cross.  This lets you know what symbol you are holding and ready to place on your schematic.  (schematic placement)</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr><w:r><w:t>Go ahead and place two of these and gates, by hovering to where you want to place them, and clicking once .  Press escape to clear your placement selection and go back to mouse mode.  (schematic placed).</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr><w:r><w:t>To see more clearly, lets zoom in.  Use the zoom tools to zoom in on your parts.   (zoom.png)</w:t></w:r><w:r><w:t xml:space="preserve">  or the zoom select tool from the verticle menu.</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr><w:r><w:t>Use the add wire tool in the verticle toolbar to add short wires to the three leads on your and gate.  Clicking on the terminal on the and gate, then clicking where you want the wire to go drops a wire segment.  Double clicking will end the wireplacement.  Press escape to end wireplacement all together and clear your tool.   (add wire tool)</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="994"/></w:tabs></w:pPr></w:p><w:p><w:r><w:t>Add IO markers.  These markers wire your circuit to the inputs and outputs of your fpga.</w:t></w:r><w:r><w:t xml:space="preserve">   </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Click the add IO button in verticle tool bar.  Click all 3 terminal leads to drop a marker.  (lead marker dropped).</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Right click to rename each port.  Name the inputs A and B and the output OUT1 (io rename)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Check your schematic for errors in Tools -&gt; Check shematic (check schematic png)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Errors show up in the conosole at the bottom.  Fix any errors.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>File -&gt; Save all</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Now we test our schematic to see if it works.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>Project -&gt; Add Source</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Select VHDL Test Bench</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Give it a name and press next  (addtb)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>The next window will ask you to select wich vhdl or schematic to test</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>You only have one (AND_schematic) so select that (tbconfirm)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Click next to view a summary (tbsummary)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>The wizard will close and you will be presented with a file called [name of file].vhd which is your vhdl code you will use to specify how to test your schematic you slected in the wizard..</w:t></w:r></w:p><w:p/><w:p/><w:p/><w:p w:rsidR="00317097" w:rsidRDefault="00421D76" w:rsidP="00421D76"/>
'@

$null = $insertRange.InsertXML($fragment)

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
